$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.966.26"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "1.891.92"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'0.7713"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.83%  "
$ws.Range("D6").Value = "'243.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.3117"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("D9").Value = "'25.72"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.01%  "
$ws.Range("D10").Value = "'0.07161"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.50%  "
$ws.Range("D11").Value = "'0.08612"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.72%  "
$ws.Range("D12").Value = "2.040.36"
$ws.Range("E12").Value = "  +7.96%  "
$ws.Range("D13").Value = "'0.7723"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.23%  "
$ws.Range("D14").Value = "'5.370"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.52%  "
$ws.Range("D15").Value = "'94.15"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.12%  "
$ws.Range("D16").Value = "'6.174"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("D17").Value = "30.164.92"
$ws.Range("E17").Value = "  +1.18%  "
$ws.Range("D18").Value = "2.352.31"
$ws.Range("E18").Value = "  +9.08%  "
$ws.Range("D19").Value = "'13.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.38%  "
$ws.Range("D20").Value = "'245.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.06%  "
$ws.Range("D21").Value = "'0.000007811"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.99%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").Value = "'8.031"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.00%  "
$ws.Range("D24").Value = "'1.002"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").Value = "'0.1642"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.21%  "
$ws.Range("D26").Value = "'9.383"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "'162.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.45%  "
$ws.Range("D28").Value = "'18.79"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.46%  "
$ws.Range("D29").Value = "'2.041"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.36%  "
$ws.Range("D30").Value = "'1.435"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.18%  "
$ws.Range("D31").Value = "'1.542"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.34%  "
$ws.Range("D32").Value = "'4.496"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.07%  "
$ws.Range("D33").Value = "'4.120"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.13%  "
$ws.Range("D34").Value = "'0.05453"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.41%  "
$ws.Range("D35").Value = "'1.241"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.03%  "
$ws.Range("D36").Value = "'0.7520"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.95%  "
$ws.Range("D37").Value = "'1.005"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.58%  "
$ws.Range("D38").Value = "'2.698"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.81%  "
$ws.Range("D39").Value = "'0.01960"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.38%  "
$ws.Range("D40").Value = "'2.786"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.21%  "
$ws.Range("D41").Value = "'0.4489"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.98%  "
$ws.Range("D42").Value = "1.111.15"
$ws.Range("E42").Value = "  -2.76%  "
$ws.Range("D43").Value = "'73.86"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.91%  "
$ws.Range("D44").Value = "'6.092"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.74%  "
$ws.Range("D45").Value = "2.237.77"
$ws.Range("E45").Value = "  +9.05%  "
$ws.Range("D46").Value = "'0.8507"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.17%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "'1.001"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "'103.96"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.33%  "
$ws.Range("D49").Value = "'1.873"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.33%  "
$ws.Range("D50").Value = "'7.612"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "'9.856"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.39%  "
